# Update the arithmetic-answer table's 100 cells to the regenerated
# values (commit "Update master to output generated at c986bee").
# Each cell is addressed directly by (row, column) so duplicate answer
# strings elsewhere in the table are never ambiguously matched, and
# setting Range.Text in place preserves the existing run formatting
# (TimeNewRoman, sz 30) already on each cell.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "97-48=49"
$t.Cell(1, 2).Range.Text = "72-33=39"
$t.Cell(1, 3).Range.Text = "85-48=37"
$t.Cell(1, 4).Range.Text = "71-68=3"
$t.Cell(1, 5).Range.Text = "90-41=49"
$t.Cell(2, 1).Range.Text = "83+8=91"
$t.Cell(2, 2).Range.Text = "44+9=53"
$t.Cell(2, 3).Range.Text = "94-9=85"
$t.Cell(2, 4).Range.Text = "95-78=17"
$t.Cell(2, 5).Range.Text = "30-5=25"
$t.Cell(3, 1).Range.Text = "34+57=91"
$t.Cell(3, 2).Range.Text = "50-9=41"
$t.Cell(3, 3).Range.Text = "16+58=74"
$t.Cell(3, 4).Range.Text = "74-56=18"
$t.Cell(3, 5).Range.Text = "48-19=29"
$t.Cell(4, 1).Range.Text = "77+4=81"
$t.Cell(4, 2).Range.Text = "61-42=19"
$t.Cell(4, 3).Range.Text = "18+58=76"
$t.Cell(4, 4).Range.Text = "63-7=56"
$t.Cell(4, 5).Range.Text = "38+57=95"
$t.Cell(5, 1).Range.Text = "16+69=85"
$t.Cell(5, 2).Range.Text = "48+13=61"
$t.Cell(5, 3).Range.Text = "53-34=19"
$t.Cell(5, 4).Range.Text = "26+37=63"
$t.Cell(5, 5).Range.Text = "22-13=9"
$t.Cell(6, 1).Range.Text = "81-35=46"
$t.Cell(6, 2).Range.Text = "71-55=16"
$t.Cell(6, 3).Range.Text = "57+35=92"
$t.Cell(6, 4).Range.Text = "95-16=79"
$t.Cell(6, 5).Range.Text = "32-26=6"
$t.Cell(7, 1).Range.Text = "65-17=48"
$t.Cell(7, 2).Range.Text = "61-49=12"
$t.Cell(7, 3).Range.Text = "9+88=97"
$t.Cell(7, 4).Range.Text = "64-19=45"
$t.Cell(7, 5).Range.Text = "50-9=41"
$t.Cell(8, 1).Range.Text = "7+5=12"
$t.Cell(8, 2).Range.Text = "37+36=73"
$t.Cell(8, 3).Range.Text = "40-12=28"
$t.Cell(8, 4).Range.Text = "31-29=2"
$t.Cell(8, 5).Range.Text = "36+59=95"
$t.Cell(9, 1).Range.Text = "13+38=51"
$t.Cell(9, 2).Range.Text = "18+16=34"
$t.Cell(9, 3).Range.Text = "27+58=85"
$t.Cell(9, 4).Range.Text = "18+16=34"
$t.Cell(9, 5).Range.Text = "91-58=33"
$t.Cell(10, 1).Range.Text = "8+69=77"
$t.Cell(10, 2).Range.Text = "83-79=4"
$t.Cell(10, 3).Range.Text = "85-78=7"
$t.Cell(10, 4).Range.Text = "48+26=74"
$t.Cell(10, 5).Range.Text = "4+89=93"
$t.Cell(11, 1).Range.Text = "16+78=94"
$t.Cell(11, 2).Range.Text = "5+39=44"
$t.Cell(11, 3).Range.Text = "5+56=61"
$t.Cell(11, 4).Range.Text = "73-49=24"
$t.Cell(11, 5).Range.Text = "47+39=86"
$t.Cell(12, 1).Range.Text = "91-45=46"
$t.Cell(12, 2).Range.Text = "28+59=87"
$t.Cell(12, 3).Range.Text = "14+79=93"
$t.Cell(12, 4).Range.Text = "81-3=78"
$t.Cell(12, 5).Range.Text = "69+8=77"
$t.Cell(13, 1).Range.Text = "21-16=5"
$t.Cell(13, 2).Range.Text = "45+16=61"
$t.Cell(13, 3).Range.Text = "42-25=17"
$t.Cell(13, 4).Range.Text = "32-19=13"
$t.Cell(13, 5).Range.Text = "67+26=93"
$t.Cell(14, 1).Range.Text = "64-17=47"
$t.Cell(14, 2).Range.Text = "47+36=83"
$t.Cell(14, 3).Range.Text = "79+8=87"
$t.Cell(14, 4).Range.Text = "75-49=26"
$t.Cell(14, 5).Range.Text = "97-19=78"
$t.Cell(15, 1).Range.Text = "73-65=8"
$t.Cell(15, 2).Range.Text = "75+18=93"
$t.Cell(15, 3).Range.Text = "92-23=69"
$t.Cell(15, 4).Range.Text = "50-23=27"
$t.Cell(15, 5).Range.Text = "53-34=19"
$t.Cell(16, 1).Range.Text = "16+38=54"
$t.Cell(16, 2).Range.Text = "49+2=51"
$t.Cell(16, 3).Range.Text = "57+4=61"
$t.Cell(16, 4).Range.Text = "26+8=34"
$t.Cell(16, 5).Range.Text = "52-33=19"
$t.Cell(17, 1).Range.Text = "90-36=54"
$t.Cell(17, 2).Range.Text = "4+79=83"
$t.Cell(17, 3).Range.Text = "5+36=41"
$t.Cell(17, 4).Range.Text = "36+7=43"
$t.Cell(17, 5).Range.Text = "70-31=39"
$t.Cell(18, 1).Range.Text = "2+29=31"
$t.Cell(18, 2).Range.Text = "19+49=68"
$t.Cell(18, 3).Range.Text = "22+59=81"
$t.Cell(18, 4).Range.Text = "35+19=54"
$t.Cell(18, 5).Range.Text = "91-89=2"
$t.Cell(19, 1).Range.Text = "29+25=54"
$t.Cell(19, 2).Range.Text = "8+66=74"
$t.Cell(19, 3).Range.Text = "49+46=95"
$t.Cell(19, 5).Range.Text = "80-33=47"
$t.Cell(20, 1).Range.Text = "45+49=94"
$t.Cell(20, 2).Range.Text = "25-16=9"
$t.Cell(20, 3).Range.Text = "47+4=51"
$t.Cell(20, 4).Range.Text = "85-59=26"
$t.Cell(20, 5).Range.Text = "24+38=62"
